# Updates Leve profit-tracking cells in Gungnir_Profits per latest market data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 57
$ws.Range("H57").Value = 21540
$ws.Range("J57").Value = 25156
$ws.Range("L57").Value = 75468
$ws.Range("N57").Value = -76466

# Row 64
$ws.Range("H64").Value = 4383.5835
$ws.Range("I64").Value = 3800
$ws.Range("J64").Value = 5200.6
$ws.Range("K64").Value = 3800
$ws.Range("L64").Value = 5200.6
$ws.Range("M64").Value = -3552
$ws.Range("N64").Value = -5696.6

# Row 67
$ws.Range("H67").Value = 4383.5835
$ws.Range("I67").Value = 3800
$ws.Range("J67").Value = 5200.6
$ws.Range("K67").Value = 3800
$ws.Range("L67").Value = 5200.6
$ws.Range("M67").Value = -2942
$ws.Range("N67").Value = -6916.6

# Row 139
$ws.Range("H139").Value = 43182.637
$ws.Range("J139").Value = 43182.637
$ws.Range("L139").Value = 43182.637
$ws.Range("N139").Value = -53462.637

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 1577.1818
$ws.Range("I63").Value = 1478.4286
$ws.Range("J63").Value = 1750
$ws.Range("K63").Value = 1478.4286
$ws.Range("L63").Value = 1750
$ws.Range("M63").Value = -792.4286
$ws.Range("N63").Value = -3122

# Row 64
$ws.Range("H64").Value = 16998
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 16998
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 16998
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -17494

# Row 66
$ws.Range("H66").Value = 1577.1818
$ws.Range("I66").Value = 1478.4286
$ws.Range("J66").Value = 1750
$ws.Range("K66").Value = 7392.143
$ws.Range("L66").Value = 8750
$ws.Range("M66").Value = -3960.143
$ws.Range("N66").Value = -15614

# Row 67
$ws.Range("H67").Value = 16998
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 16998
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 16998
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -18714

# Row 88
$ws.Range("H88").Value = 2609.923
$ws.Range("I88").Value = 3283.3333
$ws.Range("J88").Value = 2407.9
$ws.Range("K88").Value = 3283.3333
$ws.Range("L88").Value = 2407.9
$ws.Range("M88").Value = -2877.3333
$ws.Range("N88").Value = -3219.9

# Row 91
$ws.Range("H91").Value = 2609.923
$ws.Range("I91").Value = 3283.3333
$ws.Range("J91").Value = 2407.9
$ws.Range("K91").Value = 3283.3333
$ws.Range("L91").Value = 2407.9
$ws.Range("M91").Value = -1879.3333
$ws.Range("N91").Value = -5215.9

# Row 107
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

# Row 108
$ws.Range("H108").Value = 13675
$ws.Range("J108").Value = 13675
$ws.Range("L108").Value = 13675
$ws.Range("N108").Value = -21355

# Row 109
$ws.Range("H109").Value = 36943.5
$ws.Range("J109").Value = 36943.5
$ws.Range("L109").Value = 36943.5
$ws.Range("N109").Value = -39717.5

# Row 111
$ws.Range("H111").Value = 32000
$ws.Range("J111").Value = 32000
$ws.Range("L111").Value = 32000
$ws.Range("N111").Value = -40180

# Row 112
$ws.Range("H112").Value = 20354.8
$ws.Range("J112").Value = 20354.8
$ws.Range("L112").Value = 20354.8
$ws.Range("N112").Value = -23308.8

# Row 113
$ws.Range("H113").Value = 46224.5
$ws.Range("J113").Value = 46224.5
$ws.Range("L113").Value = 46224.5
$ws.Range("N113").Value = -54902.5

# Row 115
$ws.Range("H115").Value = 30000
$ws.Range("J115").Value = 30000
$ws.Range("L115").Value = 30000
$ws.Range("N115").Value = -33134

# Row 117
$ws.Range("H117").Value = 24800
$ws.Range("J117").Value = 24800
$ws.Range("L117").Value = 24800
$ws.Range("N117").Value = -33978

# Row 121
$ws.Range("H121").Value = 23254
$ws.Range("J121").Value = 23254
$ws.Range("L121").Value = 23254
$ws.Range("N121").Value = -26748

$ws = $wb.Worksheets.Item("BSM")
# Row 62
$ws.Range("H62").Value = 19111
$ws.Range("J62").Value = 19111
$ws.Range("L62").Value = 19111
$ws.Range("N62").Value = -20483

# Row 65
$ws.Range("H65").Value = 19111
$ws.Range("J65").Value = 19111
$ws.Range("L65").Value = 57333
$ws.Range("N65").Value = -64197

# Row 86
$ws.Range("H86").Value = 2328536.2
$ws.Range("I86").Value = 2600
$ws.Range("J86").Value = 3325366
$ws.Range("K86").Value = 2600
$ws.Range("L86").Value = 3325366
$ws.Range("M86").Value = -1477
$ws.Range("N86").Value = -3327612

# Row 89
$ws.Range("H89").Value = 2328536.2
$ws.Range("I89").Value = 2600
$ws.Range("J89").Value = 3325366
$ws.Range("K89").Value = 13000
$ws.Range("L89").Value = 16626830
$ws.Range("M89").Value = -7384
$ws.Range("N89").Value = -16638062

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 37038550
$ws.Range("I58").Value = 100001110
$ws.Range("J58").Value = 1753.1765
$ws.Range("K58").Value = 100001110
$ws.Range("L58").Value = 1753.1765
$ws.Range("M58").Value = -100000907
$ws.Range("N58").Value = -2159.1765

# Row 94
$ws.Range("H94").Value = 4207.2
$ws.Range("I94").Value = 2011.5
$ws.Range("K94").Value = 2011.5
$ws.Range("M94").Value = -1560.5

# Row 122
$ws.Range("H122").Value = 15625886
$ws.Range("I122").Value = 22728120
$ws.Range("J122").Value = 971.6
$ws.Range("K122").Value = 68184360
$ws.Range("L122").Value = 2914.8
$ws.Range("M122").Value = -68181910
$ws.Range("N122").Value = -7814.8

# Row 132
$ws.Range("H132").Value = 11495774
$ws.Range("I132").Value = 1256.9333
$ws.Range("J132").Value = 23811328
$ws.Range("K132").Value = 3770.7999
$ws.Range("L132").Value = 71433984
$ws.Range("M132").Value = -1240.7999
$ws.Range("N132").Value = -71439044

# Row 136
$ws.Range("H136").Value = 37038550
$ws.Range("I136").Value = 100001110
$ws.Range("J136").Value = 1753.1765
$ws.Range("K136").Value = 300003330
$ws.Range("L136").Value = 5259.529500000001
$ws.Range("M136").Value = -300000780
$ws.Range("N136").Value = -10359.5295

$ws = $wb.Worksheets.Item("GSM")
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# Row 120
$ws.Range("H120").Value = 26000
$ws.Range("J120").Value = 26000
$ws.Range("L120").Value = 26000
$ws.Range("N120").Value = -35676

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 9055.333000000001
$ws.Range("I122").Value = 12045.091
$ws.Range("K122").Value = 36135.273
$ws.Range("M122").Value = -33685.273

$ws = $wb.Worksheets.Item("WVR")
# Row 138
$ws.Range("H138").Value = 47087.918
$ws.Range("J138").Value = 47087.918
$ws.Range("L138").Value = 47087.918
$ws.Range("N138").Value = -57367.918
